$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column I: "1 - False Positive Rate" for each data row (I2:I8)
$ws.Range("I2").Formula = "=1-F2"
$ws.Range("I3").Formula = "=1-F3"
$ws.Range("I4").Formula = "=1-F4"
$ws.Range("I5").Formula = "=1-F5"
$ws.Range("I6").Formula = "=1-F6"
$ws.Range("I7").Formula = "=1-F7"
$ws.Range("I8").Formula = "=1-F8"

# Set column A width like the other bestFit columns
$ws.Columns("A").AutoFit()

# Update selection to match diff: I2:I8, active cell I2
$ws.Range("I2:I8").Select()
